# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# across the "展览" and "全部类型" sheets, matching the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value = 16320
$wsExpo.Range("G2").Value = "已售罄"
$wsExpo.Range("F3").Value = 350
$wsExpo.Range("F4").Value = 730
$wsExpo.Range("F5").Value = 252
$wsExpo.Range("F6").Value = 683
$wsExpo.Range("F7").Value = 1687
$wsExpo.Range("F8").Value = 156

# --- Sheet "全部类型" (rows 2-5 mirror "展览"; later rows shifted to 8/9/11) ---
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value = 16320
$wsAll.Range("G2").Value = "已售罄"
$wsAll.Range("F3").Value = 350
$wsAll.Range("F4").Value = 730
$wsAll.Range("F5").Value = 252
$wsAll.Range("F8").Value = 683
$wsAll.Range("F9").Value = 1687
$wsAll.Range("F11").Value = 156
